# ImplementingNonSequencedRanges/Output.xlsx was regenerated against a newer
# Aspose.Cells build (2014 -> 2016), which only changes the copyright year
# baked into the "Evaluation Only" watermark text on the "Evaluation
# Warning" sheet, and resets each sheet's page setup to use the printer's
# default paper size instead of a hard-coded "9" (Letter).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$warningSheet = $wb.Worksheets.Item("Evaluation Warning")

# Bump the Aspose.Cells evaluation watermark's copyright year from 2014 to 2016.
$warningSheet.Range("A5").Value = "Evaluation Only. Created with Aspose.Cells for .NET.Copyright 2003 - 2016 Aspose Pty Ltd."

# Clear the explicit paper size (was "9" = Letter) on both sheets so the
# workbook falls back to the default paper size.
$sheet1.PageSetup.PaperSize = 0
$warningSheet.PageSetup.PaperSize = 0
